# Apply the "Add files via upload" commit: append six new worksheets
# (Day 30 students / subjects / examinations, Day 32 SalesPerson / Company / Orders)
# after "Day 29 employee_uni", populate their data, and touch up a couple of
# pre-existing sheet selections.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the selection on "Day 29 employees" without changing which sheet
#    is active (it was not the active tab before and should not become one).
# ---------------------------------------------------------------------------
$wsEmployees = $wb.Worksheets.Item("Day 29 employees")
$wsEmployees.Range("H11").Select()

# ---------------------------------------------------------------------------
# 2. Add the six new worksheets, in order, after "Day 29 employee_uni".
#    Worksheets.Add() activates the freshly inserted sheet, which is exactly
#    what marks the previously-active "Day 29 employee_uni" tab as no longer
#    selected, and (at the very end) makes "Day 32 Orders" the active tab.
# ---------------------------------------------------------------------------
$after = $wb.Worksheets.Item($wb.Worksheets.Count)

$wsStudents     = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $after)
$wsStudents.Name = "Day 30 students"

$wsSubjects     = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsStudents)
$wsSubjects.Name = "Day 30 subjects"

$wsExaminations = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSubjects)
$wsExaminations.Name = "Day 30 examinations"

$wsSalesPerson  = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsExaminations)
$wsSalesPerson.Name = "Day 32 SalesPerson"

$wsCompany      = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsSalesPerson)
$wsCompany.Name = "Day 32 Company"

$wsOrders       = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsCompany)
$wsOrders.Name = "Day 32 Orders"

# A cell somewhere in the existing workbook that already carries the
# m/d/yyyy date style (numFmtId 14 / style index 2) -- copying its format
# instead of assigning a NumberFormat string avoids minting a brand-new
# (duplicate) number format / cell style entry.
$dateStyleSource = $wb.Worksheets.Item("Day 4").Range("D2")

# ---------------------------------------------------------------------------
# 3. "Day 30 students"
# ---------------------------------------------------------------------------
$wsStudents.Range("A1").Value = "student_id"
$wsStudents.Range("B1").Value = "student name"
$wsStudents.Range("A2").Value = 1
$wsStudents.Range("B2").Value = "Alice"
$wsStudents.Range("A3").Value = 2
$wsStudents.Range("B3").Value = "Bob"
$wsStudents.Range("A4").Value = 13
$wsStudents.Range("B4").Value = "John"
$wsStudents.Range("A5").Value = 6
$wsStudents.Range("B5").Value = "Alex"
$wsStudents.Columns.Item(1).ColumnWidth = 8.89322916666667
$wsStudents.Columns.Item(2).ColumnWidth = 11.5299479166667

# ---------------------------------------------------------------------------
# 4. "Day 30 subjects"
# ---------------------------------------------------------------------------
$wsSubjects.Range("A1").Value = "subject_name"
$wsSubjects.Range("A2").Value = "Math"
$wsSubjects.Range("A3").Value = "Physics"
$wsSubjects.Range("A4").Value = "Programming"
$wsSubjects.Columns.Item(1).ColumnWidth = 11.6197916666667
$wsSubjects.Range("D8").Select()

# ---------------------------------------------------------------------------
# 5. "Day 30 examinations"
# ---------------------------------------------------------------------------
$wsExaminations.Range("A1").Value = "student_id"
$wsExaminations.Range("B1").Value = "subject_name"

$examRows = @(
    @(1, "Math"),
    @(1, "physics"),
    @(1, "Programming"),
    @(2, "Programming"),
    @(1, "physics"),
    @(1, "Math"),
    @(13, "Math"),
    @(13, "Programming"),
    @(13, "physics"),
    @(2, "Math"),
    @(1, "Math")
)
for ($i = 0; $i -lt $examRows.Count; $i++) {
    $r = $i + 2
    $wsExaminations.Range("A$r").Value = $examRows[$i][0]
    $wsExaminations.Range("B$r").Value = $examRows[$i][1]
}
$wsExaminations.Columns.Item(1).ColumnWidth = 8.89322916666667
$wsExaminations.Columns.Item(2).ColumnWidth = 11.6197916666667
$wsExaminations.Range("J15").Select()

# ---------------------------------------------------------------------------
# 6. "Day 32 SalesPerson"
# ---------------------------------------------------------------------------
$wsSalesPerson.Range("A1").Value = "sales_id"
$wsSalesPerson.Range("B1").Value = "name"
$wsSalesPerson.Range("C1").Value = "salary"
$wsSalesPerson.Range("D1").Value = "commission_rate"
$wsSalesPerson.Range("E1").Value = "hire_date"

$wsSalesPerson.Range("A2").Value = 1
$wsSalesPerson.Range("B2").Value = "John"
$wsSalesPerson.Range("C2").Value = 100000
$wsSalesPerson.Range("D2").Value = 6
$dateStyleSource.Copy($wsSalesPerson.Range("E2"))
$wsSalesPerson.Range("E2").Value2 = 38721

$wsSalesPerson.Range("A3").Value = 2
$wsSalesPerson.Range("B3").Value = "Amy"
$wsSalesPerson.Range("C3").Value = 12000
$wsSalesPerson.Range("D3").Value = 5
$dateStyleSource.Copy($wsSalesPerson.Range("E3"))
$wsSalesPerson.Range("E3").Value2 = 40183

$wsSalesPerson.Range("A4").Value = 3
$wsSalesPerson.Range("B4").Value = "Mark"
$wsSalesPerson.Range("C4").Value = 65000
$wsSalesPerson.Range("D4").Value = 12
# hire_date for Mark is stored as a literal text string, not a real date.
$wsSalesPerson.Range("E4").Value = "'12/25/2008"
$wsSalesPerson.Range("E4").Style = "Normal"

$wsSalesPerson.Range("A5").Value = 4
$wsSalesPerson.Range("B5").Value = "Pam"
$wsSalesPerson.Range("C5").Value = 25000
$wsSalesPerson.Range("D5").Value = 25
$dateStyleSource.Copy($wsSalesPerson.Range("E5"))
$wsSalesPerson.Range("E5").Value2 = 38353

$wsSalesPerson.Range("A6").Value = 5
$wsSalesPerson.Range("B6").Value = "Alex"
$wsSalesPerson.Range("C6").Value = 5000
$wsSalesPerson.Range("D6").Value = 10
$dateStyleSource.Copy($wsSalesPerson.Range("E6"))
$wsSalesPerson.Range("E6").Value2 = 39143

$wsSalesPerson.Columns.Item(1).ColumnWidth = 6.07291666666667
$wsSalesPerson.Columns.Item(2).ColumnWidth = 4.70963541666667
$wsSalesPerson.Columns.Item(3).ColumnWidth = 5.98307291666667
$wsSalesPerson.Columns.Item(4).ColumnWidth = 14.3463541666667
$wsSalesPerson.Columns.Item(5).ColumnWidth = 9.61979166666667
$wsSalesPerson.Range("J10").Select()

# ---------------------------------------------------------------------------
# 7. "Day 32 Company"
# ---------------------------------------------------------------------------
$wsCompany.Range("A1").Value = "com_id"
$wsCompany.Range("B1").Value = "name"
$wsCompany.Range("C1").Value = "city"

$wsCompany.Range("A2").Value = 1
$wsCompany.Range("B2").Value = "RED"
$wsCompany.Range("C2").Value = "Boston"

$wsCompany.Range("A3").Value = 2
$wsCompany.Range("B3").Value = "ORANGE"
$wsCompany.Range("C3").Value = "New York"

$wsCompany.Range("A4").Value = 3
$wsCompany.Range("B4").Value = "YELLOW"
$wsCompany.Range("C4").Value = "Boston"

$wsCompany.Range("A5").Value = 4
$wsCompany.Range("B5").Value = "GREEN"
$wsCompany.Range("C5").Value = "Austin"

$wsCompany.Columns.Item(1).ColumnWidth = 6.07291666666667
$wsCompany.Columns.Item(2).ColumnWidth = 7.16666666666667
$wsCompany.Range("H14").Select()

# ---------------------------------------------------------------------------
# 8. "Day 32 Orders"
# ---------------------------------------------------------------------------
$wsOrders.Range("A1").Value = "order_id"
$wsOrders.Range("B1").Value = "order_date"
$wsOrders.Range("C1").Value = "com_id"
$wsOrders.Range("D1").Value = "sales_id"
$wsOrders.Range("E1").Value = "amount"

$wsOrders.Range("A2").Value = 1
$dateStyleSource.Copy($wsOrders.Range("B2"))
$wsOrders.Range("B2").Value2 = 41640
$wsOrders.Range("C2").Value = 3
$wsOrders.Range("D2").Value = 4
$wsOrders.Range("E2").Value = 10000

$wsOrders.Range("A3").Value = 2
$dateStyleSource.Copy($wsOrders.Range("B3"))
$wsOrders.Range("B3").Value2 = 41641
$wsOrders.Range("C3").Value = 4
$wsOrders.Range("D3").Value = 5
$wsOrders.Range("E3").Value = 5000

$wsOrders.Range("A4").Value = 3
$dateStyleSource.Copy($wsOrders.Range("B4"))
$wsOrders.Range("B4").Value2 = 41642
$wsOrders.Range("C4").Value = 1
$wsOrders.Range("D4").Value = 1
$wsOrders.Range("E4").Value = 50000

$wsOrders.Range("A5").Value = 4
$dateStyleSource.Copy($wsOrders.Range("B5"))
$wsOrders.Range("B5").Value2 = 41643
$wsOrders.Range("C5").Value = 1
$wsOrders.Range("D5").Value = 4
$wsOrders.Range("E5").Value = 25000

$wsOrders.Columns.Item(1).ColumnWidth = 7.16666666666667
$wsOrders.Columns.Item(2).ColumnWidth = 9.34635416666667
$wsOrders.Columns.Item(3).ColumnWidth = 6.07291666666667
$wsOrders.Columns.Item(4).ColumnWidth = 6.61979166666667
$wsOrders.Columns.Item(5).ColumnWidth = 6.52994791666667

# "Day 32 Orders" ends up both the last sheet and the active tab.
$wsOrders.Range("G11").Select()

Write-Output "edit complete"
